# This workbook is a quarterly financial "Income Statement" report.
# The update drops the oldest quarter column (D = "فصل دوم منتهی به 1399/06")
# and appends a new quarter column at the end (M = "فصل چهارم منتهی به 1401/12"),
# shifting all the other quarter columns one position to the left.
# In addition, the previously-reported publish date for one quarter (now in
# column I) is revised.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the oldest quarter (column D). Excel automatically shifts
#    columns E:M left into D:L, carrying their values/styles along.
$ws.Columns("D").Delete()

# 2. Recreate column M (the new last column) by copying the formatting
#    (styles/borders/number formats) of column L, the previous last column.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match column M's width to the width pattern (same as the "wide" columns,
# which is the width that column E now has after the deletion).
$ws.Columns("M").ColumnWidth = $ws.Columns("E").ColumnWidth

# 3. Fill in the new quarter's data in column M.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-06 (2)"

$ws.Range("M11").Value = 1228423
$ws.Range("M12").Value = -429191
$ws.Range("M13").Value = 799232
$ws.Range("M14").Value = -66471
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = -12848
$ws.Range("M17").Value = 719913
$ws.Range("M18").Value = -20748
$ws.Range("M19").Value = 30863
$ws.Range("M20").Value = 730028
$ws.Range("M21").Value = -59673
$ws.Range("M22").Value = 670355
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 670355
$ws.Range("M25").Value = 2235
$ws.Range("M26").Value = 300000
$ws.Range("M27").Value = 2235

# 4. Revise the publish date that had been carried over by the shift into
#    column I (previously reported as "1401-10-28 (7)", now updated).
$ws.Range("I9").Value = "1402-02-06 (9)"
